$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function ReplaceCell($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $d.Range($cell.Range.Start, $cell.Range.End)
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: row=$row col=$col old=$oldText"
    }
}

ReplaceCell 1 1 "76×54=4104" "33×47=1551"
ReplaceCell 1 2 "13×80=1040" "73×68=4964"
ReplaceCell 1 3 "46×15=690" "62×42=2604"
ReplaceCell 1 4 "61×94=5734" "79×95=7505"
ReplaceCell 1 5 "93×19=1767" "62×47=2914"

ReplaceCell 5 1 "52×11=572" "72×55=3960"
ReplaceCell 5 2 "80×36=2880" "44×28=1232"
ReplaceCell 5 3 "88×73=6424" "70×23=1610"
ReplaceCell 5 4 "14×89=1246" "43×16=688"
ReplaceCell 5 5 "44×41=1804" "42×97=4074"

ReplaceCell 10 1 "88×73=6424" "32×96=3072"
ReplaceCell 10 2 "87×75=6525" "19×26=494"
ReplaceCell 10 3 "82×33=2706" "27×83=2241"
ReplaceCell 10 4 "36×33=1188" "19×56=1064"
ReplaceCell 10 5 "35×21=735" "34×41=1394"

ReplaceCell 15 1 "55×41=2255" "52×20=1040"
ReplaceCell 15 2 "61×30=1830" "58×32=1856"
ReplaceCell 15 3 "88×90=7920" "40×39=1560"
ReplaceCell 15 4 "32×92=2944" "66×13=858"
ReplaceCell 15 5 "27×17=459" "29×92=2668"

ReplaceCell 20 1 "30×97=2910" "78×98=7644"
ReplaceCell 20 2 "28×43=1204" "40×88=3520"
ReplaceCell 20 3 "89×92=8188" "39×97=3783"
ReplaceCell 20 4 "13×20=260" "74×57=4218"
ReplaceCell 20 5 "54×19=1026" "78×28=2184"

Write-Output "Done"
